$d = $word.ActiveDocument

# --- 1. Insert the new disclaimer paragraph before the first paragraph ---
$insertionPoint = $d.Range(0, 0)
$fragment1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="360"/><w:jc w:val="both"/><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="red"/></w:rPr><w:t>Les commandes ci-dessous s</w:t></w:r><w:r><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="red"/></w:rPr><w:t xml:space="preserve">e doivent d’être exécutées </w:t></w:r><w:r><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="red"/></w:rPr><w:t>dans un terminal/invite de commandes</w:t></w:r><w:r><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="red"/></w:rPr><w:t>. Il est également possible d’utiliser le fichier « </w:t></w:r><w:r><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="red"/></w:rPr><w:t>Y-P_106-ESR-E_BackupRestore.cmd</w:t></w:r><w:r><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="red"/></w:rPr><w:t> » pour les exécuter automatiquement.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($fragment1)

# --- 2. Split the "Sauvegarde (Backup)" run into "S" + "auvegarde (Backup)" ---
$titlePara = $d.Paragraphs(2)
$titleRange = $titlePara.Range.Duplicate
$titleRange.MoveEnd(1, -1)
$fragment2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:color w:val="auto"/><w:highlight w:val="lightGray"/><w:lang w:val="fr-FR"/></w:rPr><w:t>S</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:color w:val="auto"/><w:highlight w:val="lightGray"/><w:lang w:val="fr-FR"/></w:rPr><w:t>auvegarde (Backup)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titleRange.InsertXML($fragment2)
